# Adjusting DD and DPE for EPICP according to provided dataset
$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "Variables": lower-case the three smoking-item variable names
# ------------------------------------------------------------------
$wsVars = $wb.Worksheets.Item("Variables")
$wsVars.Range("B10").Value = "zt3"
$wsVars.Range("B11").Value = "zr3"
$wsVars.Range("B12").Value = "pf3"

# ------------------------------------------------------------------
# Sheet "Categories": re-derive the category table for the
# case*_fup5 / inccanc_fup5 / vitstat5 variables (rows 39-71) to
# match the levels actually present in the dataset.
# ------------------------------------------------------------------
$wsCats = $wb.Worksheets.Item("Categories")

$rows = @(
    @("casemi_fup5", 0, "not diseased"),
    @("casemi_fup5", 1, "prevalent"),
    @("casemi_fup5", 2, "incident (verif.)"),
    @("casemi_fup5", 5, "I252 (old MI, unknown date of occurence)"),
    @("casemi_fup5", 9, "incident (not verif.)"),
    @("casestroke_fup5", 0, "not diseased"),
    @("casestroke_fup5", 1, "prevalent"),
    @("casestroke_fup5", 2, "incident (verif.)"),
    @("casestroke_fup5", 6, "I64Y old stroke, date of diagnosis unknown"),
    @("casestroke_fup5", 9, "incident (not verif.)"),
    @("caseI63_fup5", 0, "No"),
    @("caseI63_fup5", 1, "Yes"),
    @("caseI61_fup5", 0, "No"),
    @("caseI61_fup5", 1, "Yes"),
    @("casehyp_fup5", 0, "not diseased"),
    @("casehyp_fup5", 1, "prevalent"),
    @("casehyp_fup5", 2, "incident (verif.)"),
    @("casehyp_fup5", 3, "incident I15"),
    @("casehyp_fup5", 9, "incident (not verif.)"),
    @("casehf_fup5", 0, "not diseased"),
    @("casehf_fup5", 1, "prevalent"),
    @("casehf_fup5", 2, "incident (verif.)"),
    @("casehf_fup5", 9, "incident (not verif.)"),
    @("casediab_fup5", 0, "not diseased"),
    @("casediab_fup5", 1, "prevalent"),
    @("casediab_fup5", 2, "incident (verif.)"),
    @("casediab_fup5", 4, "inc. Diabetes (other types)"),
    @("casediab_fup5", 9, "incident (not verif.)"),
    @("inccanc_fup5", 0, "No"),
    @("inccanc_fup5", 1, "Yes"),
    @("vitstat5", 1, "alive"),
    @("vitstat5", 2, "dead"),
    @("vitstat5", 6, "dropped out")
)

$startRow = 39
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $entry = $rows[$i]
    $wsCats.Cells.Item($r, 1).Value = $entry[0]
    $wsCats.Cells.Item($r, 2).Value = $entry[1]
    $wsCats.Cells.Item($r, 3).Value = $entry[2]
}
